$d = $word.ActiveDocument

# Update version number 1.0 -> 1.1
$d.Content.Find.Execute("1.0", $true, $false, $false, $false, $false, $true, 1, $false, "1.1", 2)

# Update "total" -> "minimum" in the password complexity description
$d.Content.Find.Execute("8 Characters total", $true, $false, $false, $false, $false, $true, 1, $false, "8 Characters minimum", 2)
